$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (Beteckning) to know data extent.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Column C ("Förändrad") holds a date serial number for every data row (2..lastRow).
# Update it to the new date serial value 46075 (2026-02-22) for every data row.
$range = $ws.Range("C2:C$lastRow")
$range.Value2 = 46075
